# Update column F ("dSF") values for the specified rows to reflect the
# repulled data / recalculated means described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    8  = 2
    9  = -1
    21 = -2
    22 = -1
    26 = -1
    28 = 3
    31 = 1
    41 = 2
    44 = -2
    49 = 1
    52 = -1
    58 = 4
    59 = 2
    60 = 0
    64 = 0
    67 = 4
    73 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
